$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.436.86'
$ws.Range("E2").Value = '  +0.70%  '

$ws.Range("D3").Value = '1.616.50'
$ws.Range("E3").Value = '  +1.52%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.48'
$ws.Range("E5").Value = '  -0.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.498'
$ws.Range("E6").Value = '  -0.37%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.23'
$ws.Range("E10").Value = '  +1.48%  '

$ws.Range("E11").Value = '  -0.25%  '

$ws.Range("D12").Value = '1.844.47'
$ws.Range("E12").Value = '  +1.40%  '

$ws.Range("D13").Value = '1.613.19'
$ws.Range("E13").Value = '  +1.46%  '

$ws.Range("E14").Value = '  +0.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.509'
$ws.Range("E15").Value = '  +0.11%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.77'
$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '235.97'
$ws.Range("E17").Value = '  +9.63%  '

$ws.Range("D18").Value = '26.455.52'
$ws.Range("E18").Value = '  +0.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.73'
$ws.Range("E19").Value = '  +5.32%  '

$ws.Range("E20").Value = '  +0.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("E23").Value = '  +4.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.06'
$ws.Range("E24").Value = '  +0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.05'
$ws.Range("E25").Value = '  +1.58%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("E27").Value = '  +0.70%  '

$ws.Range("E28").Value = '  +0.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.51'
$ws.Range("E29").Value = '  +2.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0496'
$ws.Range("E30").Value = '  +0.96%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.16'
$ws.Range("E31").Value = '  +0.17%  '

$ws.Range("D32").Value = '1.512.30'
$ws.Range("E32").Value = '  +6.55%  '

$ws.Range("E33").Value = '  +1.67%  '

$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("E35").Value = '  +4.58%  '

$ws.Range("E36").Value = '  -0.22%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.565'
$ws.Range("E37").Value = '  -1.68%  '

$ws.Range("E38").Value = '  +0.25%  '

$ws.Range("E39").Value = '  +0.75%  '

$ws.Range("E40").Value = '  +2.18%  '

$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("E42").Value = '  +1.64%  '

$ws.Range("D43").Value = '1.757.04'
$ws.Range("E43").Value = '  +1.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.761'
$ws.Range("E44").Value = '  +0.10%  '

$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.912'
$ws.Range("E45").Value = '  -2.81%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.42'
$ws.Range("E46").Value = '  +0.93%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.97'
$ws.Range("E47").Value = '  +3.56%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").Value = '  -0.70%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.50'
$ws.Range("E49").Value = '  +1.23%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0502'
$ws.Range("E50").Value = '  +0.13%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0961'
$ws.Range("E51").Value = '  +1.10%  '
